# Update countries & provincias Spain
# Applies the 28-May-2020 04:35 data refresh to the "Pais" sheet:
#   - Bolivia's row is updated with new totals and the country ranking
#     shifts (Malasia/Marruecos/Moldavia/Ghana/Australia each drop one
#     rank, rows 61-66).
#   - Belice and Nueva Caledonia swap rank (rows 200-201).
#   - Papua Nueva Guinea / Islas Virgenes Britanicas swap rank (rows 213-214).
#   - Bonaire, San Eustaquio y Saba / San Bartolome swap rank (rows 215-216).
#   - The "Datos actualizados" timestamp moves from 04:05 to 04:35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp footer -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 04:35"

# --- Rows 61-66: Bolivia's new data shifts the ranking -----------------
$ws.Range("A61").Value = "Bolivia"
$ws.Range("B61").Value = 7768
$ws.Range("C61").Value = 632
$ws.Range("D61").Value = 689
$ws.Range("E61").Value = 6799
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 6
$ws.Range("H61").Value = 280

$ws.Range("A62").Value = "Malasia"
$ws.Range("B62").Value = 7619
$ws.Range("C62").Value = 0
$ws.Range("D62").Value = 6083
$ws.Range("E62").Value = 1421
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 115

$ws.Range("A63").Value = "Marruecos"
$ws.Range("B63").Value = 7601
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 4978
$ws.Range("E63").Value = 2421
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 202

$ws.Range("A64").Value = "Moldavia"
$ws.Range("B64").Value = 7537
$ws.Range("C64").Value = 0
$ws.Range("D64").Value = 3884
$ws.Range("E64").Value = 3379
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 274

$ws.Range("A65").Value = "Ghana"
$ws.Range("B65").Value = 7303
$ws.Range("C65").Value = 0
$ws.Range("D65").Value = 2412
$ws.Range("E65").Value = 4857
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 34

$ws.Range("A66").Value = "Australia"
$ws.Range("B66").Value = 7139
$ws.Range("C66").Value = 0
$ws.Range("D66").Value = 6566
$ws.Range("E66").Value = 470
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 103

# --- Rows 200-201: Belice / Nueva Caledonia swap rank ------------------
$ws.Range("A200").Value = "Belice"
$ws.Range("D200").Value = 16
$ws.Range("H200").Value = 2

$ws.Range("A201").Value = "Nueva Caledonia"
$ws.Range("D201").Value = 18
$ws.Range("H201").Value = 0

# --- Rows 213-214: Papua Nueva Guinea / Islas Virgenes Britanicas swap -
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1

# --- Rows 215-216: Bonaire, San Eustaquio y Saba / San Bartolome swap --
# (their numeric columns are already identical, only the names swap)
$ws.Range("A215").Value = "San Bartolome"
$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"
